$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.26126961572082763
$ws.Range("B1").Value = 0.26056370210667978
$ws.Range("A2").Value = -0.16133618743582279
$ws.Range("B2").Value = 0.15981922270004301
$ws.Range("A3").Value = -0.1101027264691119
$ws.Range("B3").Value = 0.10965690855121402
$ws.Range("A4").Value = -0.10165690867353305
$ws.Range("B4").Value = 0.10125971408806045
$ws.Range("A5").Value = -0.098259714158872136
$ws.Range("B5").Value = 0.096912831887643769
$ws.Range("A6").Value = -0.032539900222428031
$ws.Range("B6").Value = 0.032331517135226306
$ws.Range("A7").Value = -0.022331517305222537
$ws.Range("B7").Value = 0.022294981827737814
$ws.Range("A8").Value = -0.012294982000567778
$ws.Range("B8").Value = 0.012261517188800308
$ws.Range("A9").Value = -0.010261517268710385
$ws.Range("B9").Value = 0.01024230345681465
$ws.Range("A10").Value = -0.0082423035377541254
$ws.Range("B10").Value = 0.0082430038281113838
$ws.Range("A11").Value = -0.0052430039207687074
$ws.Range("B11").Value = 0.0052432815833522994
$ws.Range("A12").Value = -0.00174328168186344
$ws.Range("B12").Value = 0.0017427880000249196
$ws.Range("A13").Value = 0.0017572119022171151
$ws.Range("B13").Value = -0.001758349765785816
$ws.Range("A14").Value = 0.009758349615727191
$ws.Range("B14").Value = -0.0097717224768780753
$ws.Range("A15").Value = 0.010771722410695794
$ws.Range("B15").Value = -0.010795901396733143
$ws.Range("A16").Value = 0.01279590132067332
$ws.Range("B16").Value = -0.012918186765905837
$ws.Range("A17").Value = -0.0048008446871969923
$ws.Range("B17").Value = 0.0047988707352688564
$ws.Range("A18").Value = -0.093227514745748863
$ws.Range("B18").Value = 0.092947981562279125
$ws.Range("A19").Value = -0.053776491030842255
$ws.Range("B19").Value = 0.052946342447129702
$ws.Range("A20").Value = -0.048946342520046926
$ws.Range("B20").Value = 0.048710751553254994
$ws.Range("A21").Value = -0.004005676159747118
$ws.Range("B21").Value = 0.0039999999225273086
$ws.Range("A22").Value = -0.045716496298975784
$ws.Range("B22").Value = 0.04550237724581585
$ws.Range("A23").Value = -0.040502377329903361
$ws.Range("B23").Value = 0.040099482406427356
$ws.Range("A24").Value = -0.020099482673783697
$ws.Range("B24").Value = 0.019999999728833373
$ws.Range("A25").Value = -0.062372931740936721
$ws.Range("B25").Value = 0.062316940446446978
$ws.Range("A26").Value = -0.059816940528058637
$ws.Range("B26").Value = 0.059746803091693224
$ws.Range("A27").Value = -0.057246803176656424
$ws.Range("B27").Value = 0.056838892708016342
$ws.Range("A28").Value = -0.089011912029877571
$ws.Range("B28").Value = 0.088360133266092156
$ws.Range("A29").Value = -0.081360133432457182
$ws.Range("B29").Value = 0.081171893640935799
$ws.Range("A30").Value = -0.021171894440641825
$ws.Range("B30").Value = 0.021021753723843961
$ws.Range("A31").Value = -0.014021753903232792
$ws.Range("B31").Value = 0.01400066151233581
$ws.Range("A32").Value = -0.0040006617275185619
$ws.Range("B32").Value = 0.003999999855475167

# Match column B width to column A (stored width 15.42578125). The COM
# ColumnWidth setter snaps to a 1/6-character pixel grid, so 14.667 is the
# input that lands on the closest achievable stored width (15.5).
$ws.Columns.Item(2).ColumnWidth = 14.667
